$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row permutation map: targetRow -> sourceRow (source values, read from the
# original/before state, are written into the target row).
$map = @{
    2 = 11
    3 = 37
    4 = 16
    5 = 55
    6 = 3
    7 = 5
    8 = 7
    9 = 68
    10 = 20
    11 = 54
    12 = 19
    13 = 47
    14 = 14
    15 = 64
    16 = 40
    17 = 53
    18 = 60
    19 = 8
    20 = 76
    21 = 4
    22 = 43
    23 = 71
    24 = 35
    25 = 9
    26 = 6
    27 = 51
    28 = 42
    29 = 13
    30 = 62
    31 = 65
    32 = 22
    33 = 12
    34 = 39
    35 = 30
    36 = 80
    37 = 49
    38 = 32
    39 = 15
    40 = 74
    41 = 61
    42 = 73
    43 = 59
    44 = 67
    45 = 33
    46 = 52
    47 = 25
    48 = 50
    49 = 27
    50 = 72
    51 = 77
    52 = 38
    53 = 58
    54 = 2
    55 = 48
    56 = 46
    57 = 24
    58 = 41
    59 = 45
    60 = 31
    61 = 29
    62 = 75
    63 = 57
    64 = 63
    65 = 36
    66 = 21
    67 = 28
    68 = 70
    69 = 66
    70 = 78
    71 = 79
    72 = 34
    73 = 44
    74 = 69
    75 = 17
    76 = 23
    77 = 26
    78 = 18
    79 = 56
    80 = 10
}

# Columns that move together as a record when rows are permuted.
$cols = @(4, 10, 11, 12, 13, 14, 15, 16)  # D, J, K, L, M, N, O, P

# 1) Snapshot every source cell BEFORE any writes, so overlapping cycles in
#    the permutation do not clobber data we still need to read.
$snapshot = @{}
foreach ($r in 2..80) {
    foreach ($c in $cols) {
        $snapshot["$r`_$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# 2) Apply the permutation using the snapshot as the read source.
foreach ($targetRow in $map.Keys) {
    $sourceRow = $map[$targetRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value2 = $snapshot["$sourceRow`_$c"]
    }
}
